$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")

# Make "Repayment Schedule" the active sheet (also clears tabSelected on
# whichever sheet previously had it, e.g. "NewLoanInput").
$ws.Activate()

# Insert a new blank column before column N (14th column), shifting the
# existing N/O/P ("Late"/"Outstanding") columns right to O/P/Q.
$ws.Columns("N").Insert()

# Update the selection shown on the "Repayment Schedule" sheet.
$ws.Range("S6").Select()
